# "incluida simulacao do compras" - add a sample/simulated purchase-order
# row (row 2) below the existing header row of the COMPRAS_pedido sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# seq
$ws.Range("A2").Value = 1
# seq_cotacao
$ws.Range("B2").Value = 2
# data_pedido -> 02/05/2024
$ws.Range("C2").Value = 45414
$ws.Range("C2").NumberFormat = "dd/mm/yy"
# data_recebimento -> 10/05/2024
$ws.Range("D2").Value = 45422
$ws.Range("D2").NumberFormat = "dd/mm/yy"
# valor_pedido
$ws.Range("E2").Value = 152.5
# desconto
$ws.Range("F2").Value = 5
# situacao
$ws.Range("G2").Value = "Aguardando"

# Leave the selection back on the first cell, like the saved workbook.
[void]$ws.Range("A1").Select()
